# Changes - 19 July
# Insert a new "EventBadges" worksheet between "Contact" and
# "SubscriptionPreferences", populate it with header/data rows, and
# update the selection stored on "SubscriptionPreferences".

$wb = $excel.ActiveWorkbook

# Set the selection on SubscriptionPreferences (A1:C2) before we touch
# anything else, so its own stored cursor reflects that range.
$subSheet = $wb.Worksheets.Item("SubscriptionPreferences")
$subSheet.Range("A1:C2").Select() | Out-Null

# Insert the new worksheet right after "Contact" (i.e. right before
# "SubscriptionPreferences"), matching the target tab order:
#   Users, AdditionalInfo, Contact, EventBadges, SubscriptionPreferences
$contactSheet = $wb.Worksheets.Item("Contact")
$newSheet = $wb.Worksheets.Add($null, $contactSheet)
$newSheet.Name = "EventBadges"

# Header row (bold)
$newSheet.Range("A1").Value = "FirstName"
$newSheet.Range("B1").Value = "LastName"
$newSheet.Range("C1").Value = "CompanyName"
$newSheet.Range("A1:C1").Font.Bold = $true

# Data row
$newSheet.Range("A2").Value = "First"
$newSheet.Range("B2").Value = "Last"
$newSheet.Range("C2").Value = "ActivityCompany"

# Leave the cursor on C2 of the new (now active) sheet.
$newSheet.Range("C2").Select() | Out-Null
